$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 4396.2
$ws.Cells.Item(19, 9).Value = 989
$ws.Cells.Item(19, 10).Value = 5248
$ws.Cells.Item(19, 11).Value = 989
$ws.Cells.Item(19, 12).Value = 5248
$ws.Cells.Item(19, 13).Value = -814
$ws.Cells.Item(19, 14).Value = -5598

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 276.22223
$ws.Cells.Item(28, 9).Value = 281.88235
$ws.Cells.Item(28, 10).Value = 180
$ws.Cells.Item(28, 11).Value = 281.88235
$ws.Cells.Item(28, 12).Value = 180
$ws.Cells.Item(28, 13).Value = 203.11765
$ws.Cells.Item(28, 14).Value = -1150

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 999.5
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 14).Value = $null

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 8338098
$ws.Cells.Item(51, 10).Value = 9263998
$ws.Cells.Item(51, 12).Value = 9263998
$ws.Cells.Item(51, 14).Value = -9264966

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(88, 8).Value = 1762.125
$ws.Cells.Item(88, 10).Value = 1836.2727
$ws.Cells.Item(88, 12).Value = 1836.2727
$ws.Cells.Item(88, 14).Value = -2648.2727

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(91, 8).Value = 1762.125
$ws.Cells.Item(91, 10).Value = 1836.2727
$ws.Cells.Item(91, 12).Value = 1836.2727
$ws.Cells.Item(91, 14).Value = -4644.2727

# ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).Value = $null
$ws.Cells.Item(92, 14).Value = $null

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 1861.1025
$ws.Cells.Item(138, 9).Value = 1052.5385
$ws.Cells.Item(138, 10).Value = 3478.2307
$ws.Cells.Item(138, 11).Value = 3157.6155
$ws.Cells.Item(138, 12).Value = 10434.6921
$ws.Cells.Item(138, 13).Value = 1982.3845
$ws.Cells.Item(138, 14).Value = -20714.6921

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 35343.434
$ws.Cells.Item(32, 9).Value = 19316.947
$ws.Cells.Item(32, 11).Value = 19316.947
$ws.Cells.Item(32, 13).Value = -19029.947

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 2405
$ws.Cells.Item(63, 9).Value = 2405
$ws.Cells.Item(63, 11).Value = 2405
$ws.Cells.Item(63, 13).Value = -1719

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 2405
$ws.Cells.Item(66, 9).Value = 2405
$ws.Cells.Item(66, 11).Value = 12025
$ws.Cells.Item(66, 13).Value = -8593

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 2177.0715
$ws.Cells.Item(122, 9).Value = 1956.5834
$ws.Cells.Item(122, 11).Value = 5869.7502
$ws.Cells.Item(122, 13).Value = -3419.7502

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 4177.5
$ws.Cells.Item(86, 9).Value = 1996
$ws.Cells.Item(86, 10).Value = 4904.6665
$ws.Cells.Item(86, 11).Value = 1996
$ws.Cells.Item(86, 12).Value = 4904.6665
$ws.Cells.Item(86, 13).Value = -873
$ws.Cells.Item(86, 14).Value = -7150.6665

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 4177.5
$ws.Cells.Item(89, 9).Value = 1996
$ws.Cells.Item(89, 10).Value = 4904.6665
$ws.Cells.Item(89, 11).Value = 9980
$ws.Cells.Item(89, 12).Value = 24523.3325
$ws.Cells.Item(89, 13).Value = -4364
$ws.Cells.Item(89, 14).Value = -35755.3325

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 9443.1
$ws.Cells.Item(105, 9).Value = 10702.857
$ws.Cells.Item(105, 11).Value = 10702.857
$ws.Cells.Item(105, 13).Value = -8955.857

# BSM row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(140, 8).Value = 95347
$ws.Cells.Item(140, 10).Value = 95347
$ws.Cells.Item(140, 12).Value = 95347
$ws.Cells.Item(140, 14).Value = -105707

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1711.375
$ws.Cells.Item(16, 9).Value = 1558.6
$ws.Cells.Item(16, 11).Value = 1558.6
$ws.Cells.Item(16, 13).Value = -1271.6

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2663.5
$ws.Cells.Item(31, 9).Value = 2061
$ws.Cells.Item(31, 11).Value = 2061
$ws.Cells.Item(31, 13).Value = -1766

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2663.5
$ws.Cells.Item(34, 9).Value = 2061
$ws.Cells.Item(34, 11).Value = 2061
$ws.Cells.Item(34, 13).Value = -1859

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 4083.3333
$ws.Cells.Item(62, 10).Value = 3999.5
$ws.Cells.Item(62, 12).Value = 3999.5
$ws.Cells.Item(62, 14).Value = -5247.5

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 4083.3333
$ws.Cells.Item(65, 10).Value = 3999.5
$ws.Cells.Item(65, 12).Value = 19997.5
$ws.Cells.Item(65, 14).Value = -26237.5

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 7853.625
$ws.Cells.Item(86, 9).Value = 7388.1665
$ws.Cells.Item(86, 11).Value = 7388.1665
$ws.Cells.Item(86, 13).Value = -6265.1665

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 7853.625
$ws.Cells.Item(89, 9).Value = 7388.1665
$ws.Cells.Item(89, 11).Value = 36940.8325
$ws.Cells.Item(89, 13).Value = -31324.8325

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 1711.375
$ws.Cells.Item(113, 9).Value = 1558.6
$ws.Cells.Item(113, 11).Value = 1558.6
$ws.Cells.Item(113, 13).Value = 611.4000000000001

# CRP row 119
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 12).Value = 0
$ws.Cells.Item(119, 14).Value = $null

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2117.4783
$ws.Cells.Item(132, 9).Value = 2023.05
$ws.Cells.Item(132, 11).Value = 6069.15
$ws.Cells.Item(132, 13).Value = -3539.15

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 2411.5
$ws.Cells.Item(134, 9).Value = 2604.3333
$ws.Cells.Item(134, 11).Value = 7812.999899999999
$ws.Cells.Item(134, 13).Value = -5277.999899999999

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 500
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 500
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 1500
$ws.Cells.Item(5, 13).Value = $null
$ws.Cells.Item(5, 14).Value = -1724

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1422.5454
$ws.Cells.Item(34, 9).Value = 162.33333
$ws.Cells.Item(34, 10).Value = 2934.8
$ws.Cells.Item(34, 11).Value = 486.99999
$ws.Cells.Item(34, 12).Value = 8804.400000000001
$ws.Cells.Item(34, 13).Value = -402.99999
$ws.Cells.Item(34, 14).Value = -8972.400000000001

# CUL row 51
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(51, 8).Value = 892.2857
$ws.Cells.Item(51, 9).Value = 892.2857
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 2676.8571
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = -2216.8571
$ws.Cells.Item(51, 14).Value = $null

# CUL row 58
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(58, 8).Value = 7794.8887
$ws.Cells.Item(58, 9).Value = 3751
$ws.Cells.Item(58, 10).Value = 12849.75
$ws.Cells.Item(58, 11).Value = 11253
$ws.Cells.Item(58, 12).Value = 38549.25
$ws.Cells.Item(58, 13).Value = -11125
$ws.Cells.Item(58, 14).Value = -38805.25

# CUL row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(117, 8).Value = 4600.722
$ws.Cells.Item(117, 10).Value = 5187.6
$ws.Cells.Item(117, 12).Value = 15562.8
$ws.Cells.Item(117, 14).Value = -22446.8

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 500
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 500
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 4500
$ws.Cells.Item(135, 13).Value = $null
$ws.Cells.Item(135, 14).Value = -9570

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 330116
$ws.Cells.Item(97, 9).Value = 654232
$ws.Cells.Item(97, 11).Value = 654232
$ws.Cells.Item(97, 13).Value = -653736

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2114.9375
$ws.Cells.Item(102, 10).Value = 3735.25
$ws.Cells.Item(102, 12).Value = 3735.25
$ws.Cells.Item(102, 14).Value = -6979.25

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 21740978
$ws.Cells.Item(107, 9).Value = 1046.1666
$ws.Cells.Item(107, 10).Value = 29413896
$ws.Cells.Item(107, 11).Value = 1046.1666
$ws.Cells.Item(107, 12).Value = 29413896
$ws.Cells.Item(107, 13).Value = 873.8334
$ws.Cells.Item(107, 14).Value = -29417736

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2592.5652
$ws.Cells.Item(132, 9).Value = 2132.3684
$ws.Cells.Item(132, 11).Value = 6397.1052
$ws.Cells.Item(132, 13).Value = -3867.1052

# LTW row 51
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(51, 8).Value = 40492
$ws.Cells.Item(51, 10).Value = 40492
$ws.Cells.Item(51, 12).Value = 40492
$ws.Cells.Item(51, 14).Value = -41448

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 39061.777
$ws.Cells.Item(61, 9).Value = 1998.6
$ws.Cells.Item(61, 11).Value = 1998.6
$ws.Cells.Item(61, 13).Value = -1796.6

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 39061.777
$ws.Cells.Item(113, 9).Value = 1998.6
$ws.Cells.Item(113, 11).Value = 1998.6
$ws.Cells.Item(113, 13).Value = 171.4000000000001

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 3268
$ws.Cells.Item(132, 9).Value = 2635
$ws.Cells.Item(132, 11).Value = 7905
$ws.Cells.Item(132, 13).Value = -5375

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 6959.9
$ws.Cells.Item(122, 9).Value = 6622.1113
$ws.Cells.Item(122, 11).Value = 19866.3339
$ws.Cells.Item(122, 13).Value = -17416.3339

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 18989.475
$ws.Cells.Item(132, 9).Value = 17976.594
$ws.Cells.Item(132, 11).Value = 53929.78200000001
$ws.Cells.Item(132, 13).Value = -51399.78200000001
